$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testsheet2")

# --- Test case data: Hotels Search test case ---------------------------
# Values are written in the same order a person filling out the grid by
# hand would naturally hit the cells, so that shared-string slots land in
# the same sequence as the authored workbook.
$ws.Range("A1").Value = "TC ID/Name"
$ws.Range("A2").Value = "testingHotelsSearch"
$ws.Range("B1").Value = "Destination"
$ws.Range("C1").Value = "Checkin Date"
$ws.Range("D1").Value = "Checkout Date"
$ws.Range("B2").Value = "Grand Plaza Serviced"
$ws.Range("D2").Value = "02/02/2021"
$ws.Range("E1").Value = "Expected"

# Checkin Date (C2) is a real date serial (1/1/2021) formatted as a short
# date, not literal text like D2.
$ws.Range("C2").Value2 = 44197
$ws.Range("C2").NumberFormat = "mm-dd-yy"

# --- Column widths, best-fit to the new content -------------------------
$ws.Columns.Item(1).ColumnWidth = 16.5
$ws.Columns.Item(2).ColumnWidth = 17.16666666666667
$ws.Columns.Item(3).ColumnWidth = 10.666666666666666
$ws.Columns.Item(4).ColumnWidth = 12.0

# --- Selection / active sheet -------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("C1").Select() | Out-Null
